$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Insert two new columns (H:I) ahead of the old CAPEX column, shifting
#     CAPEX/FOM/VOM/ANR Th Eff/Carbon-intensity columns from H:M to J:O ---
$ws.Range("H1:I1").EntireColumn.Insert()

# --- Header text for the two new "equivalent total consumption" columns ---
$ws.Range("H1").Value = "Eq tot H2ElecCons (MWhe/kgh2)"
$ws.Range("I1").Value = "Eq tot H2HeatCons (MWht/kgh2)"

# --- Fill in the new per-row formulas: equivalent total electricity and
#     heat consumption for each ANR-H2 coupling, folding the ANR thermal
#     efficiency (col N) into the direct electric/heat demand (cols F/G) ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("H$r").Formula = "=F$r+G$r*N$r"
    $ws.Range("I$r").Formula = "=G$r+F$r/N$r"
}

# The engine copies the neighbouring percentage number format onto new
# formula cells that multiply by column N; strip that back off column H so
# the new cells stay "General" formatted like the rest of the row.
$ws.Range("H2:H16").ClearFormats()

# Match the width of the new columns to the existing H2ElecCons/H2HeatCons
# columns (stored width 23 == ColumnWidth 22.1666... once Excel round-trips
# the character-width<->stored-width conversion).
$ws.Range("H1:I16").ColumnWidth = 22.166666666666668

# --- View state: the Summary sheet becomes the active tab/selected sheet,
#     with the given cell selected (clears tabSelected on whichever sheet
#     had it before, matching the diff's removal from "Summary future HTSE") ---
$ws.Activate()
$ws.Range("I26").Select()
